$d = $word.ActiveDocument

# The bibliography entry "23.ed. São Paulo: Cortez, 2009." is followed, in the
# original document, by a trailing "site footer" block consisting of three
# paragraphs that must be removed entirely:
#   1) an empty paragraph
#   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3) "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github
#       pages. Original theme under Creative Commons Attribution"
# Everything else (including the empty paragraph and page-break paragraph
# that follow this block) must stay untouched.

$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 0) {
    $startParagraph = $d.Paragraphs.Item($targetIndex - 1)   # preceding empty paragraph
    $endParagraph   = $d.Paragraphs.Item($targetIndex + 1)   # the "© 2020 ..." paragraph

    $deleteRange = $d.Range($startParagraph.Range.Start, $endParagraph.Range.End)
    $deleteRange.Delete()
}
